$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 113
$ws.Range("H113").Value = 2550
$ws.Range("I113").Value = 1800
$ws.Range("J113").Value = 4800
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 4800
$ws.Range("M113").Value = 1454
$ws.Range("N113").Value = -11308

# Row 116
$ws.Range("H116").Value = 5104.5454
$ws.Range("I116").Value = 5802
$ws.Range("J116").Value = 4523.3335
$ws.Range("K116").Value = 5802
$ws.Range("L116").Value = 4523.3335
$ws.Range("M116").Value = -2360
$ws.Range("N116").Value = -11407.3335

# Row 138
$ws.Range("H138").Value = 1783.4509
$ws.Range("I138").Value = 881.125
$ws.Range("J138").Value = 3303.158
$ws.Range("K138").Value = 2643.375
$ws.Range("L138").Value = 9909.474
$ws.Range("M138").Value = 2496.625
$ws.Range("N138").Value = -20189.474

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 214258
$ws.Range("I61").Value = 1370.4828
$ws.Range("J61").Value = 557243.4399999999
$ws.Range("K61").Value = 1370.4828
$ws.Range("L61").Value = 557243.4399999999
$ws.Range("M61").Value = -1158.4828
$ws.Range("N61").Value = -557667.4399999999

# Row 136
$ws.Range("H136").Value = 214258
$ws.Range("I136").Value = 1370.4828
$ws.Range("J136").Value = 557243.4399999999
$ws.Range("K136").Value = 4111.4484
$ws.Range("L136").Value = 1671730.32
$ws.Range("M136").Value = -1561.4484
$ws.Range("N136").Value = -1676830.32

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2137.7334
$ws.Range("I86").Value = 1445.1111
$ws.Range("J86").Value = 3176.6667
$ws.Range("K86").Value = 1445.1111
$ws.Range("L86").Value = 3176.6667
$ws.Range("M86").Value = -322.1111000000001
$ws.Range("N86").Value = -5422.6667

# Row 89
$ws.Range("H89").Value = 2137.7334
$ws.Range("I89").Value = 1445.1111
$ws.Range("J89").Value = 3176.6667
$ws.Range("K89").Value = 7225.5555
$ws.Range("L89").Value = 15883.3335
$ws.Range("M89").Value = -1609.5555
$ws.Range("N89").Value = -27115.3335

# Row 97
$ws.Range("H97").Value = 10890.333
$ws.Range("I97").Value = 11000
$ws.Range("J97").Value = 10780.667
$ws.Range("K97").Value = 11000
$ws.Range("L97").Value = 10780.667
$ws.Range("M97").Value = -10009
$ws.Range("N97").Value = -12762.667

# Row 99
$ws.Range("H99").Value = 1840.2222
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 1945.25
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 1945.25
$ws.Range("M99").Value = 498
$ws.Range("N99").Value = -4941.25

# Row 102
$ws.Range("H102").Value = 2460
$ws.Range("I102").Value = 2460
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2460
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 785
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2859.353
$ws.Range("I16").Value = 1276.125
$ws.Range("J16").Value = 4266.6665
$ws.Range("K16").Value = 1276.125
$ws.Range("L16").Value = 4266.6665
$ws.Range("M16").Value = -989.125
$ws.Range("N16").Value = -4840.6665

# Row 99
$ws.Range("H99").Value = 1561.4546
$ws.Range("I99").Value = 1441.75
$ws.Range("J99").Value = 1880.6666
$ws.Range("K99").Value = 1441.75
$ws.Range("L99").Value = 1880.6666
$ws.Range("M99").Value = 56.25
$ws.Range("N99").Value = -4876.6666

# Row 113
$ws.Range("H113").Value = 2859.353
$ws.Range("I113").Value = 1276.125
$ws.Range("J113").Value = 4266.6665
$ws.Range("K113").Value = 1276.125
$ws.Range("L113").Value = 4266.6665
$ws.Range("M113").Value = 893.875
$ws.Range("N113").Value = -8606.666499999999

# Row 126
$ws.Range("H126").Value = 1561.4546
$ws.Range("I126").Value = 1441.75
$ws.Range("J126").Value = 1880.6666
$ws.Range("K126").Value = 4325.25
$ws.Range("L126").Value = 5641.9998
$ws.Range("M126").Value = -1855.25
$ws.Range("N126").Value = -10581.9998

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 1443.2307
$ws.Range("I107").Value = 1987.5
$ws.Range("J107").Value = 1201.3334
$ws.Range("K107").Value = 5962.5
$ws.Range("L107").Value = 3604.0002
$ws.Range("M107").Value = -4042.5
$ws.Range("N107").Value = -7444.0002

# Row 133
$ws.Range("H133").Value = 560871.3
$ws.Range("I133").Value = 4918.5713
$ws.Range("J133").Value = 914659.4399999999
$ws.Range("K133").Value = 14755.7139
$ws.Range("L133").Value = 2743978.32
$ws.Range("M133").Value = -9695.713899999999
$ws.Range("N133").Value = -2754098.32

# Row 138
$ws.Range("H138").Value = 2769.1667
$ws.Range("I138").Value = 2575.7144
$ws.Range("J138").Value = 3040
$ws.Range("K138").Value = 7727.1432
$ws.Range("L138").Value = 9120
$ws.Range("M138").Value = -2587.1432
$ws.Range("N138").Value = -19400

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 4317.273
$ws.Range("I102").Value = 5354
$ws.Range("J102").Value = 3453.3333
$ws.Range("K102").Value = 5354
$ws.Range("L102").Value = 3453.3333
$ws.Range("M102").Value = -3732
$ws.Range("N102").Value = -6697.3333

# Row 122
$ws.Range("H122").Value = 1138.0714
$ws.Range("I122").Value = 1104.1111
$ws.Range("J122").Value = 1199.2
$ws.Range("K122").Value = 3312.3333
$ws.Range("L122").Value = 3597.6
$ws.Range("M122").Value = -862.3333000000002
$ws.Range("N122").Value = -8497.6

# Row 126
$ws.Range("H126").Value = 1287.9
$ws.Range("I126").Value = 1046
$ws.Range("J126").Value = 1852.3334
$ws.Range("K126").Value = 3138
$ws.Range("L126").Value = 5557.0002
$ws.Range("M126").Value = -668
$ws.Range("N126").Value = -10497.0002

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 1948.7693
$ws.Range("I40").Value = 1396.0555
$ws.Range("J40").Value = 3192.375
$ws.Range("K40").Value = 1396.0555
$ws.Range("L40").Value = 3192.375
$ws.Range("M40").Value = -1260.0555
$ws.Range("N40").Value = -3464.375

# Row 122
$ws.Range("H122").Value = 31361.914
$ws.Range("I122").Value = 49905.145
$ws.Range("J122").Value = 3547.0715
$ws.Range("K122").Value = 149715.435
$ws.Range("L122").Value = 10641.2145
$ws.Range("M122").Value = -147265.435
$ws.Range("N122").Value = -15541.2145

# Row 132
$ws.Range("H132").Value = 13899.421
$ws.Range("I132").Value = 5358.1377
$ws.Range("J132").Value = 41421.332
$ws.Range("K132").Value = 16074.4131
$ws.Range("L132").Value = 124263.996
$ws.Range("M132").Value = -13544.4131
$ws.Range("N132").Value = -129323.996

$ws = $wb.Worksheets.Item("WVR")
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

# Row 122
$ws.Range("H122").Value = 12501394
$ws.Range("I122").Value = 20001100
$ws.Range("J122").Value = 1884.1666
$ws.Range("K122").Value = 60003300
$ws.Range("L122").Value = 5652.4998
$ws.Range("M122").Value = -60000850
$ws.Range("N122").Value = -10552.4998
